# CoLDP parent-child schema template: add typeMaterial / typeLocality (and
# related new fields) to the workbook.
#
# Summary of changes (see commit message "Adds type material and locality"):
#   - Description sheet: new "format" column inserted after "category"
#   - Name sheet: new "cultivarEpithet", "appendedPhrase", "original",
#     "typeMaterial" and "typeLocality" columns inserted at various points
#   - Reference sheet: new "remarks" column appended at the end
#   - Synonym sheet: new "referenceID" column inserted before "remarks"

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Description sheet: A=taxonID B=category C=format(NEW) D=description
#                     E=language F=referenceID
# ---------------------------------------------------------------------
$wsDescription = $wb.Worksheets.Item("Description")
$wsDescription.Columns.Item(3).Insert()
$wsDescription.Range("C1").Value = "format"

# ---------------------------------------------------------------------
# Name sheet: insert new columns working right-to-left so earlier column
# letters stay valid while we work.
#   ... G=infraspecificEpithet H=cultivarEpithet(NEW) I=appendedPhrase(NEW)
#   J=publishedInID K=publishedInPage L=original(NEW) M=code N=status
#   O=typeMaterial(NEW) P=typeLocality(NEW) Q=link R=remarks
# ---------------------------------------------------------------------
$wsName = $wb.Worksheets.Item("Name")

# Insert two columns before the original column L ("link") -> typeMaterial, typeLocality
$wsName.Range("L1:M1").EntireColumn.Insert()
$wsName.Range("L1").Value = "typeMaterial"
$wsName.Range("M1").Value = "typeLocality"

# Insert one column before the original column J ("code") -> original
$wsName.Columns.Item(10).Insert()
$wsName.Range("J1").Value = "original"

# Insert two columns before the original column H ("publishedInID") -> cultivarEpithet, appendedPhrase
$wsName.Range("H1:I1").EntireColumn.Insert()
$wsName.Range("H1").Value = "cultivarEpithet"
$wsName.Range("I1").Value = "appendedPhrase"

# ---------------------------------------------------------------------
# Reference sheet: append new "remarks" column after the last used column.
# ---------------------------------------------------------------------
$wsReference = $wb.Worksheets.Item("Reference")
$wsReference.Range("I1").Value = "remarks"

# ---------------------------------------------------------------------
# Synonym sheet: A=taxonID B=nameID C=status D=referenceID(NEW) E=remarks
# ---------------------------------------------------------------------
$wsSynonym = $wb.Worksheets.Item("Synonym")
$wsSynonym.Columns.Item(4).Insert()
$wsSynonym.Range("D1").Value = "referenceID"

# ---------------------------------------------------------------------
# Update row-1 selections to mirror the header-row selection left behind
# by the editing session, and make sure the first sheet ("Description")
# ends up as the active tab (matches the saved workbook view state).
# ---------------------------------------------------------------------
$wsSynonym.Rows.Item(1).Select()
$wsReference.Rows.Item(1).Select()
$wsName.Range("H1").Select()
$wsDescription.Rows.Item(1).Select()
